$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '29.811.09'
$ws.Cells.Item(2, 5).Value = '  -0.56%  '

$ws.Cells.Item(3, 4).Value = '1.890.94'
$ws.Cells.Item(3, 5).Value = '  -0.06%  '

$ws.Cells.Item(4, 5).Value = '  +0.15%  '

$ws.Cells.Item(5, 4).Formula = '="0.7921"'
$ws.Cells.Item(5, 4).Copy()
$ws.Cells.Item(5, 4).PasteSpecial(-4163)
$ws.Cells.Item(5, 5).Value = '  -4.59%  '

$ws.Cells.Item(6, 4).Formula = '="243.11"'
$ws.Cells.Item(6, 4).Copy()
$ws.Cells.Item(6, 4).PasteSpecial(-4163)
$ws.Cells.Item(6, 5).Value = '  +0.45%  '

$ws.Cells.Item(7, 4).Formula = '="1.000"'
$ws.Cells.Item(7, 4).Copy()
$ws.Cells.Item(7, 4).PasteSpecial(-4163)
$ws.Cells.Item(7, 5).Value = '  +0.10%  '

$ws.Cells.Item(8, 4).Formula = '="0.3148"'
$ws.Cells.Item(8, 4).Copy()
$ws.Cells.Item(8, 4).PasteSpecial(-4163)
$ws.Cells.Item(8, 5).Value = '  -2.96%  '

$ws.Cells.Item(9, 4).Formula = '="25.29"'
$ws.Cells.Item(9, 4).Copy()
$ws.Cells.Item(9, 4).PasteSpecial(-4163)
$ws.Cells.Item(9, 5).Value = '  -4.33%  '

$ws.Cells.Item(10, 4).Formula = '="0.07079"'
$ws.Cells.Item(10, 4).Copy()
$ws.Cells.Item(10, 4).PasteSpecial(-4163)
$ws.Cells.Item(10, 5).Value = '  +0.71%  '

$ws.Cells.Item(11, 4).Formula = '="0.08073"'
$ws.Cells.Item(11, 4).Copy()
$ws.Cells.Item(11, 4).PasteSpecial(-4163)
$ws.Cells.Item(11, 5).Value = '  +0.29%  '

$ws.Cells.Item(12, 4).Formula = '="0.7662"'
$ws.Cells.Item(12, 4).Copy()
$ws.Cells.Item(12, 4).PasteSpecial(-4163)
$ws.Cells.Item(12, 5).Value = '  +2.22%  '

$ws.Cells.Item(13, 4).Value = '1.901.02'
$ws.Cells.Item(13, 5).Value = '  +0.54%  '

$ws.Cells.Item(14, 4).Formula = '="5.385"'
$ws.Cells.Item(14, 4).Copy()
$ws.Cells.Item(14, 4).PasteSpecial(-4163)
$ws.Cells.Item(14, 5).Value = '  +2.99%  '

$ws.Cells.Item(15, 4).Formula = '="92.23"'
$ws.Cells.Item(15, 4).Copy()
$ws.Cells.Item(15, 4).PasteSpecial(-4163)
$ws.Cells.Item(15, 5).Value = '  -0.17%  '

$ws.Cells.Item(16, 4).Value = '29.823.39'
$ws.Cells.Item(16, 5).Value = '  -0.47%  '

$ws.Cells.Item(17, 4).Formula = '="5.987"'
$ws.Cells.Item(17, 4).Copy()
$ws.Cells.Item(17, 4).PasteSpecial(-4163)
$ws.Cells.Item(17, 5).Value = '  +2.08%  '

$ws.Cells.Item(18, 4).Formula = '="13.81"'
$ws.Cells.Item(18, 4).Copy()
$ws.Cells.Item(18, 4).PasteSpecial(-4163)
$ws.Cells.Item(18, 5).Value = '  -1.79%  '

$ws.Cells.Item(19, 4).Formula = '="243.41"'
$ws.Cells.Item(19, 4).Copy()
$ws.Cells.Item(19, 4).PasteSpecial(-4163)
$ws.Cells.Item(19, 5).Value = '  -1.06%  '

$ws.Cells.Item(20, 4).Formula = '="0.000007681"'
$ws.Cells.Item(20, 4).Copy()
$ws.Cells.Item(20, 4).PasteSpecial(-4163)
$ws.Cells.Item(20, 5).Value = '  -1.08%  '

$ws.Cells.Item(21, 4).Formula = '="8.246"'
$ws.Cells.Item(21, 4).Copy()
$ws.Cells.Item(21, 4).PasteSpecial(-4163)
$ws.Cells.Item(21, 5).Value = '  +18.57%  '

$ws.Cells.Item(22, 4).Formula = '="1.001"'
$ws.Cells.Item(22, 4).Copy()
$ws.Cells.Item(22, 4).PasteSpecial(-4163)
$ws.Cells.Item(22, 5).Value = '  +0.16%  '

$ws.Cells.Item(23, 4).Value = '2.140.19'
$ws.Cells.Item(23, 5).Value = '  +0.16%  '

$ws.Cells.Item(24, 5).Value = '  +0.18%  '

$ws.Cells.Item(25, 4).Formula = '="0.1648"'
$ws.Cells.Item(25, 4).Copy()
$ws.Cells.Item(25, 4).PasteSpecial(-4163)
$ws.Cells.Item(25, 5).Value = '  +3.53%  '

$ws.Cells.Item(26, 4).Formula = '="9.339"'
$ws.Cells.Item(26, 4).Copy()
$ws.Cells.Item(26, 4).PasteSpecial(-4163)
$ws.Cells.Item(26, 5).Value = '  +1.39%  '

$ws.Cells.Item(27, 4).Formula = '="165.72"'
$ws.Cells.Item(27, 4).Copy()
$ws.Cells.Item(27, 4).PasteSpecial(-4163)
$ws.Cells.Item(27, 5).Value = '  +0.04%  '

$ws.Cells.Item(28, 4).Formula = '="18.68"'
$ws.Cells.Item(28, 4).Copy()
$ws.Cells.Item(28, 4).PasteSpecial(-4163)
$ws.Cells.Item(28, 5).Value = '  -0.91%  '

$ws.Cells.Item(29, 4).Formula = '="2.041"'
$ws.Cells.Item(29, 4).Copy()
$ws.Cells.Item(29, 4).PasteSpecial(-4163)
$ws.Cells.Item(29, 5).Value = '  -2.58%  '

$ws.Cells.Item(30, 4).Formula = '="1.406"'
$ws.Cells.Item(30, 4).Copy()
$ws.Cells.Item(30, 4).PasteSpecial(-4163)
$ws.Cells.Item(30, 5).Value = '  +3.29%  '

$ws.Cells.Item(31, 5).Value = '  +1.65%  '

$ws.Cells.Item(32, 4).Formula = '="4.437"'
$ws.Cells.Item(32, 4).Copy()
$ws.Cells.Item(32, 4).PasteSpecial(-4163)
$ws.Cells.Item(32, 5).Value = '  +3.88%  '

$ws.Cells.Item(33, 4).Formula = '="0.05605"'
$ws.Cells.Item(33, 4).Copy()
$ws.Cells.Item(33, 4).PasteSpecial(-4163)
$ws.Cells.Item(33, 5).Value = '  -2.03%  '

$ws.Cells.Item(34, 4).Formula = '="4.036"'
$ws.Cells.Item(34, 4).Copy()
$ws.Cells.Item(34, 4).PasteSpecial(-4163)
$ws.Cells.Item(34, 5).Value = '  -0.72%  '

$ws.Cells.Item(35, 5).Value = '  -1.02%  '

$ws.Cells.Item(36, 4).Formula = '="0.7386"'
$ws.Cells.Item(36, 4).Copy()
$ws.Cells.Item(36, 4).PasteSpecial(-4163)
$ws.Cells.Item(36, 5).Value = '  +1.73%  '

$ws.Cells.Item(37, 4).Formula = '="1.001"'
$ws.Cells.Item(37, 4).Copy()
$ws.Cells.Item(37, 4).PasteSpecial(-4163)
$ws.Cells.Item(37, 5).Value = '  +0.29%  '

$ws.Cells.Item(38, 4).Formula = '="2.635"'
$ws.Cells.Item(38, 4).Copy()
$ws.Cells.Item(38, 4).PasteSpecial(-4163)
$ws.Cells.Item(38, 5).Value = '  -2.88%  '

$ws.Cells.Item(39, 4).Formula = '="0.01905"'
$ws.Cells.Item(39, 4).Copy()
$ws.Cells.Item(39, 4).PasteSpecial(-4163)
$ws.Cells.Item(39, 5).Value = '  -0.75%  '

$ws.Cells.Item(40, 4).Formula = '="2.780"'
$ws.Cells.Item(40, 4).Copy()
$ws.Cells.Item(40, 4).PasteSpecial(-4163)
$ws.Cells.Item(40, 5).Value = '  +0.60%  '

$ws.Cells.Item(41, 4).Formula = '="0.4400"'
$ws.Cells.Item(41, 4).Copy()
$ws.Cells.Item(41, 4).PasteSpecial(-4163)
$ws.Cells.Item(41, 5).Value = '  -0.79%  '

$ws.Cells.Item(42, 4).Formula = '="72.25"'
$ws.Cells.Item(42, 4).Copy()
$ws.Cells.Item(42, 4).PasteSpecial(-4163)
$ws.Cells.Item(42, 5).Value = '  +0.53%  '

$ws.Cells.Item(43, 2).Value = 'Maker'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(43, 4).Value = '1.058.87'
$ws.Cells.Item(43, 5).Value = '  +6.80%  '

$ws.Cells.Item(44, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(44, 4).Formula = '="0.8525"'
$ws.Cells.Item(44, 4).Copy()
$ws.Cells.Item(44, 4).PasteSpecial(-4163)
$ws.Cells.Item(44, 5).Value = '  +0.67%  '

$ws.Cells.Item(45, 2).Value = 'FraxShare'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(45, 4).Formula = '="5.795"'
$ws.Cells.Item(45, 4).Copy()
$ws.Cells.Item(45, 4).PasteSpecial(-4163)
$ws.Cells.Item(45, 5).Value = '  -2.06%  '

$ws.Cells.Item(46, 2).Value = 'PaxDollar'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(46, 4).Formula = '="1.000"'
$ws.Cells.Item(46, 4).Copy()
$ws.Cells.Item(46, 4).PasteSpecial(-4163)
$ws.Cells.Item(46, 5).Value = '  +0.10%  '

$ws.Cells.Item(47, 4).Formula = '="103.20"'
$ws.Cells.Item(47, 4).Copy()
$ws.Cells.Item(47, 4).PasteSpecial(-4163)
$ws.Cells.Item(47, 5).Value = '  +2.16%  '

$ws.Cells.Item(48, 4).Formula = '="9.991"'
$ws.Cells.Item(48, 4).Copy()
$ws.Cells.Item(48, 4).PasteSpecial(-4163)
$ws.Cells.Item(48, 5).Value = '  +2.58%  '

$ws.Cells.Item(49, 4).Formula = '="1.867"'
$ws.Cells.Item(49, 4).Copy()
$ws.Cells.Item(49, 4).PasteSpecial(-4163)
$ws.Cells.Item(49, 5).Value = '  -0.32%  '

$ws.Cells.Item(50, 4).Formula = '="7.405"'
$ws.Cells.Item(50, 4).Copy()
$ws.Cells.Item(50, 4).PasteSpecial(-4163)
$ws.Cells.Item(50, 5).Value = '  -2.13%  '

$ws.Cells.Item(51, 4).Value = '2.034.98'
$ws.Cells.Item(51, 5).Value = '  -0.06%  '

$excel.CutCopyMode = $false